$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Section 1 (rows 1-16): Uplink budget
# ---------------------------------------------------------------------------

# Pointing Loss C7: 0.2 -> 0.3
$ws.Range("C7").Value = 0.3

# New reference/notes table entries (antenna gain scratch data) added
# alongside rows 9-17 in columns F:H. Order matches the shared-string
# insertion order of the authored workbook: monopole, downlink, uplink,
# dipole, turnstile, VHF Downlink.
$ws.Range("G9").Value = "monopole"

$ws.Range("F11").Value = "downlink"
$ws.Range("G11").Value = 6.8
$ws.Range("H11").Value = 63.95

$ws.Range("F10").Value = "uplink"
$ws.Range("G10").Value = 4.7
$ws.Range("H10").Value = 72.489999999999995

$ws.Range("G12").Value = "dipole"

$ws.Range("F13").Value = "uplink"
$ws.Range("G13").Value = 0.3
$ws.Range("H13").Value = 81.290000000000006

$ws.Range("F14").Value = "downlink"
$ws.Range("G14").Value = 0.1
$ws.Range("H14").Value = 64.05

$ws.Range("G15").Value = "turnstile"

$ws.Range("F16").Value = "uplink"
$ws.Range("G16").Value = 0.3
$ws.Range("H16").Value = 81.290000000000006

# Row 13/15 formulas reference the external (unresolvable in this sandbox)
# defined name "k" ([1]Input!$C$3 == -228.6). Re-enter the formulas with the
# literal substituted in place of the external name so the dependent chain
# recalculates to the correct cached results (C15, which has no external
# reference of its own, then recalculates automatically off the corrected
# C13).
$ws.Range("C13").Formula = "=C11-C7-(-228.6)+C12"

# ---------------------------------------------------------------------------
# Row 17: Downlink Frequency input - no longer a formula, literal value, and
# relabelled from UHF to VHF.
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 145
$ws.Range("D17").Value = "VHF Downlink"

$ws.Range("F17").Value = "downlink"
$ws.Range("G17").Value = 0.1
$ws.Range("H17").Value = 64.05

# ---------------------------------------------------------------------------
# Section 2 (rows 18-31): Downlink budget
# ---------------------------------------------------------------------------
$ws.Rows(20).RowHeight = 30

# Pointing Loss C22: 0.2 -> 6.8
$ws.Range("C22").Value = 6.8

# Polarization Loss C23: was formula "=loss" (external name), now a literal.
$ws.Range("C23").Value = 0.1

# C25 formula also references the external defined name "c_"
# ([1]Input!$C$4 == 299792400). Re-enter with the literal substituted so the
# dependent chain (C26, C28 via C13-style fix below, C30) recalculates
# correctly.
$ws.Range("C25").Formula = "=22+20*LOG10((C24*1000)/(299792400/(C17*10^6)))"

# C28 formula also references the external defined name "k".
$ws.Range("C28").Formula = "=C26-C22-(-228.6)+C27"

# ---------------------------------------------------------------------------
# Cosmetic: widen the new notes column and restore the selection like the
# authored workbook.
# ---------------------------------------------------------------------------
$ws.Columns("G").AutoFit()
$ws.Range("D17").Select()
